# ProductHistorical_TestData.xlsx -- "code changes for product historical features"
#
# Renumbers the test-data fixtures from the "58/59" generation to "60", and
# normalizes "Web Data 58 Child" -> "Child Web Data 60" phrasing across all
# sheets, then leaves the CreateOrder sheet as the active tab/selection.

$wb = $excel.ActiveWorkbook

# --- LoginSignup -----------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("LoginSignup")
$wsLogin.Range("C1").Value = "Web Data 60"
$wsLogin.Range("D1").Value = "Web Data 60"
$wsLogin.Range("N1").Value = "Child Web Data 60"
$wsLogin.Range("U1").Value = "Successfully created Child Web Data 60. You can now login with the username admin after your password is set. Password reset link is sent to your email."

# Touch the far-right cell of row 1 (mirrors the extra AMJ1 cell / widened
# used-range picked up in the saved workbook).
$loginTouch = $wsLogin.Cells.Item(1, 1024)
$loginTouch.NumberFormat = "General"

# --- AddCurrency -------------------------------------------------------------
$wsCurrency = $wb.Worksheets.Item("AddCurrency")
$wsCurrency.Range("C1").Value = "Web Data 60"
$wsCurrency.Range("D1").Value = "Child Web Data 60"
$wsCurrency.Range("F1").Value = "Working as admin Child Web Data 60 X"

# --- AddProductCategory1 ----------------------------------------------------
$wsCategory = $wb.Worksheets.Item("AddProductCategory1")
$wsCategory.Range("C1").Value = "Web Data 60"
$wsCategory.Range("E1").Value = "Child Web Data 60"

# --- AddCustomer -------------------------------------------------------------
$wsCustomer = $wb.Worksheets.Item("AddCustomer")
$wsCustomer.Range("C1").Value = "Web Data 60"
$wsCustomer.Range("C2").Value = "Child Web Data 60"

# --- CreateOrder -------------------------------------------------------------
$wsOrder = $wb.Worksheets.Item("CreateOrder")
$wsOrder.Range("C1").Value = "Web Data 60"
$wsOrder.Range("C2").Value = "Child Web Data 60"

# Update the selections left in each sheet so the saved view matches.
$wsLogin.Range("E1").Select()
$wsCurrency.Range("C1").Select()
$wsCategory.Range("E1").Select()
$wsCustomer.Range("D1").Select()
$wsOrder.Range("C2").Select()

# CreateOrder ends up the active/selected tab.
$wsOrder.Activate()
